# Sprint Backlog Burndown - "Update the sprint backlog to reflect today's effort"
#
# The burndown tracker (Sheet1) records, per task, how many points were
# worked off in Week 1 (col D), Week 2 (col E) and Week 3 (col F). Today's
# effort bumped the numbers for three "Login" tasks (rows 3-5) and closed
# out the remaining point on the "Register -> Implement register process"
# task (row 7). Row 29's SUM() totals (and the burndown chart that reads
# them) recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 - Login / Create credentials data base
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0

# Row 4 - Login / Implement credentials validation process
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1

# Row 5 - Login / Create register page
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1

# Row 7 - Register / Implement register process
$ws.Range("E7").Value = 2

# Recalculate so the Estimate Totals row (29) and the burndown chart
# pick up the new figures.
$excel.CalculateFull()

# Match the reviewer's on-screen state when they saved: zoomed to 100%
# with C5 selected.
$excel.ActiveWindow.Zoom = 100
$ws.Range("C5").Select()
